$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record as row 6 (weekly refresh of the daily price
# series) -- this pushes every existing row 6..27 down by one (to 7..28),
# growing the used range from A1:R27 to A1:R28.
$ws.Rows("6:6").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 44687
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 100112001
$ws.Range("G6").Value = "Berenjena"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 9000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 9500
$ws.Range("N6").Value = "`$/caja 60 unidades"
$ws.Range("O6").Value = "Región de Arica y Parinacota"
$ws.Range("P6").Value = 158
$ws.Range("Q6").Value = 60
$ws.Range("R6").Value = "Hortaliza"
